# Update cryptocurrency price/volume data in Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '51.861.87'
$ws.Range('E2').Value = '  -0.44%  '
$ws.Range('D3').Value = '2.768.67'
$ws.Range('E3').Value = '  -2.13%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'355.41"
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('D6').Value = "'108.94"
$ws.Range('E6').Value = '  -4.03%  '
$ws.Range('E7').Value = '  +2.75%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = "'0.587"
$ws.Range('E9').Value = '  -3.15%  '
$ws.Range('D10').Value = "'40.01"
$ws.Range('E10').Value = '  -4.85%  '
$ws.Range('D11').Value = "'0.0850"
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('D12').Value = "'0.131"
$ws.Range('E12').Value = '  +0.72%  '
$ws.Range('D13').Value = "'19.34"
$ws.Range('E13').Value = '  -3.91%  '
$ws.Range('E14').Value = '  -2.47%  '
$ws.Range('D15').Value = '3.209.44'
$ws.Range('E15').Value = '  -1.56%  '
$ws.Range('D16').Value = '2.777.84'
$ws.Range('E16').Value = '  -1.94%  '
$ws.Range('D17').Value = "'0.921"
$ws.Range('E17').Value = '  +2.89%  '
$ws.Range('D18').Value = '51.736.84'
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('E19').Value = '  +0.94%  '
$ws.Range('E20').Value = '  -1.53%  '
$ws.Range('D21').Value = "'12.99"
$ws.Range('E21').Value = '  -5.67%  '
$ws.Range('D22').Value = '0.0₃0972'
$ws.Range('E22').Value = '  -2.23%  '
$ws.Range('D23').Value = "'272.82"
$ws.Range('E23').Value = '  +0.95%  '
$ws.Range('D24').Value = "'69.50"
$ws.Range('E24').Value = '  -0.27%  '
$ws.Range('E25').Value = '  -2.60%  '
$ws.Range('D26').Value = "'26.44"
$ws.Range('E26').Value = '  -0.89%  '
$ws.Range('D27').Value = "'1.00"
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('D28').Value = "'10.08"
$ws.Range('E28').Value = '  -1.90%  '
$ws.Range('E29').Value = '  -1.05%  '
$ws.Range('D30').Value = "'0.143"
$ws.Range('E30').Value = '  +1.98%  '
$ws.Range('D31').Value = "'51.31"
$ws.Range('E31').Value = '  +1.14%  '
$ws.Range('D32').Value = "'0.0459"
$ws.Range('E32').Value = '  +3.51%  '
$ws.Range('D33').Value = "'33.75"
$ws.Range('E33').Value = '  -0.32%  '
$ws.Range('E34').Value = '  -3.24%  '
$ws.Range('D35').Value = "'5.36"
$ws.Range('E35').Value = '  +9.81%  '
$ws.Range('D36').Value = "'0.0834"
$ws.Range('E36').Value = '  +0.30%  '
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('E38').Value = '  -0.58%  '
$ws.Range('D39').Value = "'18.11"
$ws.Range('E39').Value = '  -1.76%  '
$ws.Range('E40').Value = '  -5.12%  '
$ws.Range('E41').Value = '  -0.85%  '
$ws.Range('D42').Value = "'2.50"
$ws.Range('E42').Value = '  -3.52%  '
$ws.Range('D43').Value = "'123.34"
$ws.Range('E43').Value = '  -3.55%  '
$ws.Range('D44').Value = "'2.25"
$ws.Range('E44').Value = '  -2.09%  '
$ws.Range('D45').Value = "'21.70"
$ws.Range('E45').Value = '  -7.71%  '
$ws.Range('D46').Value = '2.060.06'
$ws.Range('E46').Value = '  +0.77%  '
$ws.Range('D47').Value = "'3.24"
$ws.Range('E47').Value = '  -3.48%  '
$ws.Range('E48').Value = '  -1.25%  '
$ws.Range('D49').Value = "'5.67"
$ws.Range('E49').Value = '  -0.35%  '
$ws.Range('D50').Value = "'0.920"
$ws.Range('E50').Value = '  -4.13%  '
$ws.Range('D51').Value = "'8.91"
$ws.Range('E51').Value = '  +0.02%  '
